$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'goalkeeper pants youth soccer'
$ws.Cells.Item(2, 1).Value = 'compression pants men cheap'
$ws.Cells.Item(3, 1).Value = 'leggings baseball'
$ws.Cells.Item(4, 1).Value = 'indoor knee pads'
$ws.Cells.Item(5, 1).Value = 'yoga knee pads 2 pack'
$ws.Cells.Item(6, 1).Value = 'compression shorts for basketball'
$ws.Cells.Item(7, 1).Value = 'knee pads gym'
$ws.Cells.Item(8, 1).Value = 'boys basketball leggings youth'
$ws.Cells.Item(9, 1).Value = 'knee sweat pants'
$ws.Cells.Item(10, 1).Value = 'knee pads for sports'
$ws.Cells.Item(11, 1).Value = 'black softball pants'
$ws.Cells.Item(12, 1).Value = 'arthritis hope knee compression sleeve'
$ws.Cells.Item(13, 1).Value = 'honeycomb tights'
$ws.Cells.Item(14, 1).Value = 'mens sheer pants'
$ws.Cells.Item(15, 1).Value = 'mens basketball knee sleeves'
$ws.Cells.Item(16, 1).Value = 'boys baseball pants short'
$ws.Cells.Item(17, 1).Value = 'spandex work pants men'
$ws.Cells.Item(18, 1).Value = 'baseball pants long'
$ws.Cells.Item(19, 1).Value = 'capri pants for men'
$ws.Cells.Item(20, 1).Value = 'boys basketball tights'
$ws.Cells.Item(21, 1).Value = 'tight pants'
$ws.Cells.Item(22, 1).Value = 'knee length shorts for men'
$ws.Cells.Item(23, 1).Value = 'softball shorts for men'
$ws.Cells.Item(24, 1).Value = 'baseball catchers hand pad'
$ws.Cells.Item(25, 1).Value = 'cold pad knee'
$ws.Cells.Item(26, 1).Value = 'football pants mens with pads'
$ws.Cells.Item(27, 1).Value = 'yoga positions chart'
$ws.Cells.Item(28, 1).Value = 'hex leg sleeve youth'
$ws.Cells.Item(29, 1).Value = 'knees pads yoga'
$ws.Cells.Item(30, 1).Value = 'work pants for men with knee pad'
$ws.Cells.Item(31, 1).Value = 'tight basketball shorts'
$ws.Cells.Item(32, 1).Value = 'adult hockey pants'
$ws.Cells.Item(33, 1).Value = 'youth sliding shorts baseball'
$ws.Cells.Item(34, 1).Value = 'sheer pants for men'
$ws.Cells.Item(35, 1).Value = 'paintball shorts'
$ws.Cells.Item(36, 1).Value = 'lacrosse padded shorts'
$ws.Cells.Item(37, 1).Value = 'knee pads for paintball'
$ws.Cells.Item(38, 1).Value = 'mens big and tall basketball pants'
$ws.Cells.Item(39, 1).Value = 'mens shorts below the knee'
$ws.Cells.Item(40, 1).Value = 'basketball tights youth'
$ws.Cells.Item(41, 1).Value = 'girls sliding shorts softball'
$ws.Cells.Item(42, 1).Value = 'youth compression pants boys'
$ws.Cells.Item(43, 1).Value = 'boys compression running pants'
$ws.Cells.Item(44, 1).Value = 'compression yoga tights'
$ws.Cells.Item(45, 1).Value = 'sliding shorts youth girls softball'
$ws.Cells.Item(46, 1).Value = 'basketball lot'
$ws.Cells.Item(47, 1).Value = 'paintball leg pads'
$ws.Cells.Item(48, 1).Value = 'adult black football pants'
$ws.Cells.Item(49, 1).Value = 'knee hockey pads'
$ws.Cells.Item(50, 1).Value = 'knee sleeve basketball'
$ws.Cells.Item(51, 1).Value = 'knee pads work pants'
$ws.Cells.Item(52, 1).Value = 'mens knee pads'
$ws.Cells.Item(53, 1).Value = 'padded baseball sliding shorts'
$ws.Cells.Item(54, 1).Value = 'soccer tights for men'
$ws.Cells.Item(55, 1).Value = 'boys tight pants'
$ws.Cells.Item(56, 1).Value = 'black leggings for men'
$ws.Cells.Item(57, 1).Value = 'men athletic compression pants'
$ws.Cells.Item(58, 1).Value = 'youth 5 pad girdle'
$ws.Cells.Item(59, 1).Value = 'poc knee pads'
$ws.Cells.Item(60, 1).Value = 'men baseball pants black'
$ws.Cells.Item(61, 1).Value = 'youth small compression pants'
$ws.Cells.Item(62, 1).Value = 'hex knee pads compression leg sleeve'
$ws.Cells.Item(63, 1).Value = 'hockey compression pants'
$ws.Cells.Item(64, 1).Value = 'softball mens'
$ws.Cells.Item(65, 1).Value = 'x compression pants'
$ws.Cells.Item(66, 1).Value = 'men compression legging'
$ws.Cells.Item(67, 1).Value = 'wrestling knee pads'
$ws.Cells.Item(68, 1).Value = 'boys tights for sports youth'
$ws.Cells.Item(69, 1).Value = 'hockey pants youth'
$ws.Cells.Item(70, 1).Value = 'pant sport men'
$ws.Cells.Item(71, 1).Value = 'girls softball sliding pants'
$ws.Cells.Item(72, 1).Value = 'boys small compression pants'
$ws.Cells.Item(73, 1).Value = 'compression spandex for men'
$ws.Cells.Item(74, 1).Value = 'knee pads mtb'
$ws.Cells.Item(75, 1).Value = 'sport leggings men'
$ws.Cells.Item(76, 1).Value = 'mens sliding shorts'
$ws.Cells.Item(77, 1).Value = 'kneeling on the promises'
$ws.Cells.Item(78, 1).Value = 'tall baseball pants mens'
$ws.Cells.Item(79, 1).Value = 'yoga knee pad thick'
$ws.Cells.Item(80, 1).Value = 'men gym leggings'
$ws.Cells.Item(81, 1).Value = 'basketball knee sleeves for men'
$ws.Cells.Item(82, 1).Value = 'basketball knee sleeves with pads'
$ws.Cells.Item(83, 1).Value = 'youth athletic tights'
$ws.Cells.Item(84, 1).Value = 'baseball pants adults'
$ws.Cells.Item(85, 1).Value = 'tights for soccer'
$ws.Cells.Item(86, 1).Value = 'sliding shorts'
$ws.Cells.Item(87, 1).Value = 'mens sport leggings'
$ws.Cells.Item(88, 1).Value = '6 inch basketball'
$ws.Cells.Item(89, 1).Value = 'football waist pads'
$ws.Cells.Item(90, 1).Value = 'mens work pants knee pads'
$ws.Cells.Item(91, 1).Value = 'knee sleeve lacrosse'
$ws.Cells.Item(92, 1).Value = 'men sliding shorts'
$ws.Cells.Item(93, 1).Value = 'padded leggings'
$ws.Cells.Item(94, 1).Value = 'yoga pads'
$ws.Cells.Item(95, 1).Value = 'mens cold tights'
$ws.Cells.Item(96, 1).Value = 'boy tights youth'
$ws.Cells.Item(97, 1).Value = 'basketball cycling'
$ws.Cells.Item(98, 1).Value = 'basketball knee sleeves youth'
$ws.Cells.Item(99, 1).Value = 'soccer compression pants'
$ws.Cells.Item(100, 1).Value = 'sliding pads'
